$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row for "Walfaanaa Magarsaa" (row 3) and the row for
# "Lalisee Magarsaa" (originally row 7). Delete from bottom first so
# the earlier row index is not affected by the shift.
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(3).Delete()

# Update the selection to match the saved workbook view state.
$ws.Range("C12").Select()
